# Daily scrape refresh for the AIESEC Global Talent opportunities sheet.
# - Rows 2-13 are overwritten with the newly scraped opportunities.
# - Rows 14-16 (stale opportunities) are removed entirely.
# - The "PREMIUM = Yes" highlight no longer applies to any remaining row,
#   so the old yellow highlight is cleared from the rows that used to carry it.
# - A few column widths are nudged to fit the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Opportunity IDs are digit strings scraped verbatim from the source site
# (same as every other text column here) -- keep them as text instead of
# letting them silently become numbers.
$ws.Range("A2:A13").NumberFormat = "@"

# New data: Row, OpportunityID, Link, Title, Country, Premium, Applicants, Duration, Organization
$data = @(
    ,@(2, '1330719', 'https://aiesec.org/opportunity/global-talent/1330719', 'Junior HR', 'Islamabad, Pakistan', 'No', '2 applicants', '9 - 12 Weeks', 'Enigma Software Solutions')
    ,@(3, '1330718', 'https://aiesec.org/opportunity/global-talent/1330718', 'Social Media Executive', 'Islamabad, Pakistan', 'No', '0 applicants', '9 - 12 Weeks', 'Roshan Studios')
    ,@(4, '1330717', 'https://aiesec.org/opportunity/global-talent/1330717', 'Business Development Executive', 'Rawalpindi, Pakistan', 'No', '1 applicant', '9 - 12 Weeks', 'GrowUp Tech Solution')
    ,@(5, '1330658', 'https://aiesec.org/opportunity/global-talent/1330658', 'Content Creator', 'Lahore, Pakistan', 'No', '0 applicants', '9 - 12 Weeks', 'Talent Trellis')
    ,@(6, '1330657', 'https://aiesec.org/opportunity/global-talent/1330657', 'Business Development Manager', 'Lahore, Pakistan', 'No', '1 applicant', '9 - 12 Weeks', 'Talent Trellis')
    ,@(7, '1330655', 'https://aiesec.org/opportunity/global-talent/1330655', 'Software Engineer Ruby on Rails/React/Flutter/AI/ML', 'Lahore, Pakistan', 'No', '2 applicants', '9 - 12 Weeks', 'Arkhitech')
    ,@(8, '1330643', 'https://aiesec.org/opportunity/global-talent/1330643', 'Market Data Management', 'Heerlen, Netherlands', 'No', '5 applicants', '6 - 18 Months', 'APG Heerlen')
    ,@(9, '1330625', 'https://aiesec.org/opportunity/global-talent/1330625', 'Marketing Executive', 'Islamabad, Pakistan', 'No', '0 applicants', '9 - 12 Weeks', 'GrowUp Tech Solution')
    ,@(10, '1328274', 'https://aiesec.org/opportunity/global-talent/1328274', 'Web Designer', 'Tunis, Tunisia', 'No', '27 applicants', '9 - 12 Weeks', 'La fabrique')
    ,@(11, '1328032', 'https://aiesec.org/opportunity/global-talent/1328032', 'IT Analyst', '2620 Ramada, Portugal', 'No', '61 applicants', '9 - 12 Weeks', 'Pegadamotriz')
    ,@(12, '1325464', 'https://aiesec.org/opportunity/global-talent/1325464', 'Accelerate Romania|Account Manager for Foreign Markets', 'Bucharest, Romania', 'No', '67 applicants', '9 - 12 Weeks', 'Azuvioo')
    ,@(13, '1318464', 'https://aiesec.org/opportunity/global-talent/1318464', 'Summer Internship Program on AI & ML', 'Ghaziabad, Uttar Pradesh, India', 'No', '49 applicants', '9 - 12 Weeks', 'KIET group of institutions')
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}

# Rows 2,3,5,6 used to be PREMIUM="Yes" with a yellow highlight on column E;
# they are now "No", so drop the old highlight style back to the default.
$ws.Range("E2").Style = "Normal"
$ws.Range("E3").Style = "Normal"
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Style = "Normal"

# The scrape now only has 12 opportunities (rows 2-13); drop the trailing
# rows that fell out of the feed.
$ws.Rows("14:16").Delete()

# A few columns were resized to better fit the refreshed content.
$ws.Columns("C").ColumnWidth = 57 - 5/6
$ws.Columns("D").ColumnWidth = 34 - 5/6
$ws.Columns("F").ColumnWidth = 16 - 5/6
$ws.Columns("H").ColumnWidth = 29 - 5/6
